# Synchronize the BOM: swap the LED (LNJ237W82RA) and resistor (ERJ-3GEYJ102V)
# rows on the BOM worksheet, remove the stray blank spacer row, and tidy up
# the row heights / view settings that Excel rewrote on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Swap the data held in rows 14 and 15 (columns A:J) ---------------------
$row14 = $ws.Range("A14:J14").Value2
$row15 = $ws.Range("A15:J15").Value2

$ws.Range("A14:J14").Value2 = $row15
$ws.Range("A15:J15").Value2 = $row14

# --- Swap the hyperlinks that go with that data so they keep following it --
$ws.Hyperlinks.Item($ws.Range("J14")).Delete()
$ws.Hyperlinks.Item($ws.Range("J15")).Delete()

$ws.Hyperlinks.Add($ws.Range("J14"), "http://www.digikey.com/product-detail/en/LNJ237W82RA/LNJ237W82RACT-ND/2349015", "", "", "LNJ237W82RACT-ND")
$ws.Hyperlinks.Add($ws.Range("J15"), "http://www.digikey.com/product-detail/en/ERJ-3GEYJ102V/P1.0KGCT-ND/134874", "", "", "P1.0KGCT-ND")

# --- Delete the now-unneeded blank spacer row (old row 16) -----------------
$ws.Rows.Item(16).Delete()

# --- Reset row heights on the data rows so Excel writes its default dy -----
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(16).RowHeight = 15

# --- Update the sheet view: zoom to 100%, select the whole row 15 ----------
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("A15:XFD15").Select()
$ws.Range("A15").Activate()
